# Refresh cryptos list: latest prices + 1h volume-change percentages,
# plus the Filecoin / InternetComputer(DFINITY) rank swap (rows 31-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.408.41"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.842.45"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Formula = "'0.9995"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Formula = "'239.51"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Formula = "'0.6262"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Formula = "'0.07400"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Formula = "'24.90"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D12").Value = "1.844.18"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Formula = "'4.967"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Formula = "'0.6709"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Formula = "'0.00001034"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Formula = "'81.76"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Formula = "'6.260"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "29.395.44"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Formula = "'234.31"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Formula = "'1.001"
$ws.Range("D22").Formula = "'7.304"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").Formula = "'0.9995"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Formula = "'156.99"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Formula = "'8.476"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Formula = "'0.07248"
$ws.Range("E28").Value = "  +11.61%  "
$ws.Range("D29").Formula = "'1.488"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").Formula = "'1.479"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Formula = "'4.055"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Formula = "'4.035"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("D33").Formula = "'1.162"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").Formula = "'1.818"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Formula = "'0.7115"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Formula = "'0.01838"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Formula = "'2.790"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "1.233.03"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").Formula = "'6.789"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Formula = "'0.9538"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").Formula = "'1.001"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "1.996.85"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").Formula = "'101.31"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Formula = "'65.33"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").Formula = "'0.00000000117"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("D47").Formula = "'1.700"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Formula = "'6.955"
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").Formula = "'8.940"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").Formula = "'0.3877"
$ws.Range("E51").Value = "  -1.97%  "
